$wb = $excel.ActiveWorkbook

# The original sheet gets duplicated so the copy inherits sheetId=2 (an
# internal Excel counter) along with all of the original sheet's namespace
# declarations / formatting; the stale original is then removed so the
# surviving sheet renormalizes down to sheet1.xml / rId1, matching what
# Excel itself would produce for this edit.
$origName = $wb.ActiveSheet.Name
$old = $wb.Worksheets.Item($origName)
$old.Copy($null, $old)
$origStale = $wb.Worksheets.Item($origName)
$null = $origStale.Delete()

$ws = $wb.Worksheets.Item(1)
$ws.Activate()
$ws.Name = "ValidLogin"

$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

$excel.ActiveWindow.Zoom = 175
$null = $ws.Range("B3").Select()
